$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet held totals for years 2000/2002/2005/2007/2010/2012/2015/2017
# (rows 2-9). The edit drops the four oldest years (2000, 2002, 2005, 2007)
# and adds a new 2020 row, leaving 2010/2012/2015/2017/2020 (rows 2-6).

# Deleting rows 2:5 removes 2000/2002/2005/2007 and shifts the
# 2010/2012/2015/2017 rows up into rows 2-5 (dimension becomes A1:S6
# automatically).
$ws.Rows("2:5").Delete()

# Row 6 is now free - this is where the new 2020 data goes. Copy the
# year-label formatting (bold/centered/bordered) from A2 so the new label
# matches the existing year cells, then set its text.
$ws.Range("A2").Copy()
$ws.Range("A6").PasteSpecial(-4122)
$ws.Range("A6").Value = "2020年"

$ws.Range("B6").Value = 164181553.680022
$ws.Range("C6").Value = 3790999303.43253
$ws.Range("D6").Value = 1331683250.27097
# E6 (化学工业总产出金额) has no reported figure for 2020 - leave it blank,
# matching the other blank cells already present for 2017.
$ws.Range("E6").Value = "'"
$ws.Range("E6").Style = "Normal"
$ws.Range("F6").Value = 2860304981.41469
$ws.Range("G6").Value = 26990278045.4323
# H6 (房地产业...) - no figure, blank.
$ws.Range("H6").Value = "'"
$ws.Range("H6").Style = "Normal"
# I6 (批发零售贸易...) - no figure, blank.
$ws.Range("I6").Value = "'"
$ws.Range("I6").Style = "Normal"
$ws.Range("J6").Value = 3668921660.99905
$ws.Range("K6").Value = 1922677593.56453
$ws.Range("L6").Value = 821452123.004553
$ws.Range("M6").Value = 757386926.81411
# N6 (运输仓储邮政...) - no figure, blank.
$ws.Range("N6").Value = "'"
$ws.Range("N6").Style = "Normal"
$ws.Range("O6").Value = 556533552.939994
$ws.Range("P6").Value = 1727661637.78313
# Q6 (金融业总产出金额) - no figure, blank.
$ws.Range("Q6").Value = "'"
$ws.Range("Q6").Style = "Normal"
$ws.Range("R6").Value = 733031886.5509059
$ws.Range("S6").Value = 1251362798.17111
